# 7.10 Fixed Some Bugs
# Wrap the three question lines in green color tags (used for inner-voice / hint text)
# and grow their row heights to fit the now-longer text, then move the active
# selection to B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dialogue text for rows 2-4 (B column), keeping the same wording
# but wrapping it with the <color=#00CC00>(...)</color> rich-text tag used by
# the in-game dialogue renderer.
$ws.Range("B2").Value = " <color=#00CC00>(Among these people, who has the greatest ability to destroy the bridge?)</color>"
$ws.Range("B3").Value = " <color=#00CC00>(Who has the strongest motive to do so?)</color>"
$ws.Range("B4").Value = " <color=#00CC00>(Let’s review the character profiles we have so far.)</color>"

# Rows 3 and 4 now need extra height to display the longer text (row 2 was
# already sized for two lines).
$ws.Rows.Item(3).RowHeight = 34
$ws.Rows.Item(4).RowHeight = 34

# Move/refresh the active selection as recorded in the saved view state.
$ws.Range("B14").Select()
